$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kayitlar")

# A new record (Kayıt No 5, 2025-07-10) was inserted above the existing
# "Kayıt No 2" row, pushing that row from row 2 down to row 3.
$ws.Rows.Item(2).Insert()

# All columns in this sheet store plain text (even numeric-looking values
# and dates), so prefix with a literal-text quote to avoid Excel's
# automatic number/date conversion.
$ws.Range("A2").Value = "'5"
$ws.Range("B2").Value = "'2025-07-10"
$ws.Range("C2").Value = "'Merkez"
$ws.Range("D2").Value = "'2"
$ws.Range("E2").Value = "'2"
$ws.Range("F2").Value = "'Tevhid"
$ws.Range("G2").Value = "'Gökhan ELGÜL, Göktan ELGÜL"
